# Updated cryptos list with GitHub Actions - refresh prices / 1h volumes.
# Values that look like plain decimal numbers (single '.') are written with a
# leading apostrophe so Excel keeps them as literal text (matching the
# source data's inlineStr cells) instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.270.68"
$ws.Range("E2").Value = "  -0.13%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.593.85"
$ws.Range("E3").Value = "  +0.16%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'213.02"
$ws.Range("E5").Value = "  +0.44%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.499"
$ws.Range("E6").Value = "  -0.44%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.20%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.42%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'18.98"
$ws.Range("E10").Value = "  -1.93%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.37%  "

# Row 12 - Wrapped liquid staked Ether 2.0
$ws.Range("D12").Value = "1.817.61"
$ws.Range("E12").Value = "  +0.13%  "

# Row 13
$ws.Range("D13").Value = "1.617.14"
$ws.Range("E13").Value = "  +1.74%  "

# Row 14
$ws.Range("E14").Value = "  -0.98%  "

# Row 15
$ws.Range("D15").Value = "'0.510"
$ws.Range("E15").Value = "  -2.18%  "

# Row 16
$ws.Range("D16").Value = "'63.82"
$ws.Range("E16").Value = "  -0.98%  "

# Row 17
$ws.Range("D17").Value = "26.252.45"
$ws.Range("E17").Value = "  -0.26%  "

# Row 18
$ws.Range("E18").Value = "  -1.38%  "

# Row 19
$ws.Range("D19").Value = "'215.10"
$ws.Range("E19").Value = "  +1.49%  "

# Row 20
$ws.Range("E20").Value = "  -1.82%  "

# Row 21
$ws.Range("E21").Value = "  -0.03%  "

# Row 22
$ws.Range("E22").Value = "  -0.04%  "

# Row 23
$ws.Range("D23").Value = "'9.04"
$ws.Range("E23").Value = "  +0.43%  "

# Row 24
$ws.Range("E24").Value = "  -2.57%  "

# Row 25
$ws.Range("D25").Value = "'144.84"
$ws.Range("E25").Value = "  -0.26%  "

# Row 26
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("E27").Value = "  -1.15%  "

# Row 28
$ws.Range("E28").Value = "  +0.66%  "

# Row 29
$ws.Range("E29").Value = "  -0.67%  "

# Row 30
$ws.Range("E30").Value = "  -2.15%  "

# Row 31
$ws.Range("E31").Value = "  +0.49%  "

# Row 33
$ws.Range("D33").Value = "1.416.19"
$ws.Range("E33").Value = "  +6.05%  "

# Row 34
$ws.Range("E34").Value = "  +0.01%  "

# Row 35
$ws.Range("E35").Value = "  -0.71%  "

# Row 36
$ws.Range("E36").Value = "  -1.55%  "

# Row 37
$ws.Range("D37").Value = "'0.572"
$ws.Range("E37").Value = "  -5.16%  "

# Row 38
$ws.Range("E38").Value = "  -1.01%  "

# Row 39
$ws.Range("D39").Value = "'0.822"
$ws.Range("E39").Value = "  +0.56%  "

# Row 40
$ws.Range("E40").Value = "  -0.18%  "

# Row 41
$ws.Range("E41").Value = "  -0.01%  "

# Row 42 / 43 - WEMIXToken and MXToken swapped order + refreshed values
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.16"
$ws.Range("E42").Value = "  +0.78%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'0.934"
$ws.Range("E43").Value = "  -11.96%  "

# Row 44
$ws.Range("D44").Value = "'0.761"
$ws.Range("E44").Value = "  -0.33%  "

# Row 45
$ws.Range("D45").Value = "1.729.58"
$ws.Range("E45").Value = "  +0.12%  "

# Row 46
$ws.Range("D46").Value = "'60.90"
$ws.Range("E46").Value = "  -1.56%  "

# Row 47
$ws.Range("D47").Value = "'87.48"
$ws.Range("E47").Value = "  -0.48%  "

# Row 48
$ws.Range("E48").Value = "  -1.69%  "

# Row 49
$ws.Range("E49").Value = "  -0.42%  "

# Row 50
$ws.Range("D50").Value = "'0.0954"
$ws.Range("E50").Value = "  -2.74%  "

# Row 51
$ws.Range("E51").Value = "  +0.04%  "
